# "added social media links + more"
#
# The seven original bullet paragraphs are reworked into an eight-bullet
# list: two leading bullets are kept, "Random placement" becomes a
# (grammar-flagged) "shrink images" bullet, "Give credit" moves up, the
# "_GoBack" bookmark now sits inside a new "Clean code/add comments"
# bullet, and three brand-new bullets are appended at the end.
#
# Rebuild the whole body in one shot via Range.InsertXML so the final
# paragraph/run/bookmark structure exactly matches the target, rather
# than trying to patch seven separate paragraphs in place.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$bodyXml = @"
<w:p xmlns:w="$wNs"><w:r><w:t xml:space="preserve">Top panel </w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Add photo picture, random</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:proofErr w:type="gramStart"/><w:r><w:t>shrink</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> images</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Give credit</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Clean cod</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>e/add comments</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Check browser compatibility</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Add anchors for links</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>Add div when it resizes</w:t></w:r></w:p>
"@

$d.Content.InsertXML($bodyXml)
